$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G35").Value = "u"
$ws.Range("G36").Value = "u"
$ws.Range("H41").Value = "u"
$ws.Range("I41").Value = "u"
$ws.Range("J41").Value = "u"
$ws.Range("H42").Value = "u"
$ws.Range("I42").Value = "u"
$ws.Range("H43").Value = "u"
$ws.Range("G44").Value = "u"
$ws.Range("G45").Value = "u"
$ws.Range("G46").Value = "u"

$ws.Range("H43").Select() | Out-Null
